# Regenerate merged AHB files
#
# For the "group header" rows in the 27003 sheet (rows that introduce a new
# segment, e.g. row 58, 62, 69, ...), re-apply the un-highlighted formatting
# (style index 2 / 3 for column B) that the rest of the already-regenerated
# rows (e.g. row 2) use, and clear out the stale "ÄNDERUNG" marker that used
# to live in column L for every data row in the 58-105 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 58-105 are the newly (re)generated segment block. Column L ("Änderung")
# no longer carries the "ÄNDERUNG" shared-string flag for any row in the
# block - clear its contents and restore the plain/no-highlight style (the
# same style already used by column L on e.g. row 2).
$ws.Range("L2").Copy()
$ws.Range("L58:L105").PasteSpecial(-4122)
$ws.Range("L58:L105").ClearContents()
$excel.CutCopyMode = $false

# The rows that start a new segment group get their whole row restyled from
# the "highlighted" look (style 5 / 7) to the plain look (style 2 / 3 / 4)
# already used elsewhere in the sheet (e.g. row 2), without touching the
# cell values.
$headerRows = @(58, 62, 69, 74, 77, 82, 86, 89, 94, 98, 103)

$ws.Range("A2:V2").Copy()
foreach ($r in $headerRows) {
    $target = $ws.Range("A" + $r + ":V" + $r)
    $target.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
